$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sample data: First Name, Last Name, Username, Password, Images (filename)
$data = @(
    @("Steve",   "Sadhwani",    "stevo",           "stevo",       ""),
    @("Paul",    "Hafer",       "pauleatstoomuch", "paul",        ""),
    @("Zack",    "Wakeley",     "zackman40",       "zack",        ""),
    @("Mike",    "Fuentes",     "mikey",           "mike",        "mf1.jpg,mf2.jpg,mf3.jpg,mf4.jpg"),
    @("Bob",     "Barker",      "bobbyb",          "bob",         "bb1.jpg,bb2.jpg,bb3.jpg,bb4.jpg,bb5.jpg"),
    @("Ralph",   "Tindell",     "rtindell",        "ralph",       ""),
    @("Patrick", "Finnegan",    "patrickf",        "patrick",     "pf1.jpg,pf2.jpg,pf3.jpg,pf4.jpg,pf5.jpg"),
    @("Mary",    "Francis",     "maryf",           "mary",        ""),
    @("Joey",    "Merchant",    "joeym",           "joey",        "jm1.jpg,jm2.jpg,jm3.jpg,jm4.jpg,jm5.jpg"),
    @("Alfred",  "Pennyworth",  "batman",          "masterbruce", "ap1.jpg,ap2.jpg,ap3.jpg,ap4.jpg,ap5.jpg")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    if ($r[4] -ne "") {
        $ws.Cells.Item($row, 5).Value = $r[4]
    }
    $row++
}

# Update the saved view: scrolled to show column E, selection on E7
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("E7").Select()
